# Auto commit at 2025-11-15 8:35:14.12
# Append the next day's (2025-11-14) charging-station data for both
# stations (四方坪站 / 高岭站) to the bottom of the daily data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at row 27; new rows go to 28 and 29.
$lastRow = 27
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Copy formatting (styles/number formats only) from the previous two rows
# down into the new rows so they match the rest of the table exactly,
# without disturbing the shared number-format definitions.
$srcRange = $ws.Range("A" + ($lastRow - 1) + ":F" + $lastRow)
$srcRange.Copy()
$dstRange = $ws.Range("A" + $newRow1 + ":F" + $newRow2)
$dstRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Date value as the raw serial number (2025-11-14), matching the other
# date cells and avoiding an implicit reformat that a DateTime .Value
# assignment would trigger on the already-pasted cell format.
$newDateSerial = 45975

# Row 28: 四方坪站 (station 1)
$ws.Cells.Item($newRow1, 1).Value = $newDateSerial
$ws.Cells.Item($newRow1, 2).Value = "四方坪站"
$ws.Cells.Item($newRow1, 3).Value = 9179.99
$ws.Cells.Item($newRow1, 4).Value = 8140.57
$ws.Cells.Item($newRow1, 5).Value = 3068.72
$ws.Cells.Item($newRow1, 6).Value = 380

# Row 29: 高岭站 (station 2)
$ws.Cells.Item($newRow2, 1).Value = $newDateSerial
$ws.Cells.Item($newRow2, 2).Value = "高岭站"
$ws.Cells.Item($newRow2, 3).Value = 4566.62
$ws.Cells.Item($newRow2, 4).Value = 4001.11
$ws.Cells.Item($newRow2, 5).Value = 1188.81
$ws.Cells.Item($newRow2, 6).Value = 159

$ws.Range("L34").Select() | Out-Null
